$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.026.14"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.303.00"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'100.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("D7").Value = "'0.501"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").Value = "'36.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.16%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'18.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.64%  "
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "2.662.25"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "2.314.77"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "42.914.88"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "'12.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.52%  "
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'68.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("E23").Value = "  +14.56%  "
$ws.Range("D24").Value = "'236.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "'24.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").Value = "'2.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.38%  "
$ws.Range("D29").Value = "'169.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "'34.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "'9.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'5.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "'17.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("D35").Value = "'4.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").Value = "'0.0695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "1.992.26"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").Value = "'0.0290"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("E44").Value = "  -7.04%  "
$ws.Range("D45").Value = "'10.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'17.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'56.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.528.74"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -1.52%  "
